# Applies the change: add a new Task entry (D19) describing InferSent,
# matching the style/format of the other "D column" task notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$newText = "nferSent – Supervised Learning of Sentence Embeddings/Attention mechanism etc.: https://yashuseth.wordpress.com/2018/08/06/infersent-supervised-learning-of-sentence-embeddings/"

# Write the new value into D19 (new row at the bottom of the task list)
$ws.Range("D19").Value = $newText

# Match formatting of the neighboring D-column note cells (D17/D18): wrap text
$ws.Range("D19").WrapText = $true
$ws.Range("D19").Style = $ws.Range("D18").Style

# Update the view so the newly added row/selection is visible, matching the
# workbook's saved view state after this edit.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("D20").Select()
